# Update the handback-status report timestamps ("Generate Report for Handback").
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" (Overview) / "Correspond Handoff Datetime" (de-de)
# both reference the same original timestamp value.
$wsOverview.Range("G2").Value = "2016-09-02 03:15:12"
$wsDeDe.Range("H2").Value     = "2016-09-02 03:15:12"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-02 03:15:00"
$wsZhCn.Range("K2").Value = "2016-09-02 03:15:31"

# de-de: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-09-02 03:15:38"
